$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New dataset: sending clusters ECs, FAPs, sCs (Angpt1 -> Itga5) against target clusters ECs, FAPs, M2, sCs
$senders = @("ECs", "FAPs", "sCs")
$targets = @("ECs", "FAPs", "M2", "sCs")

$values = @{
  "ECs|ECs"   = @(1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.004971127078137211, 0.004971127078137211, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 2.917670120532222, 26.25903108479, 0.001430246438017095, 0.001430246438017095)
  "ECs|FAPs"  = @(1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.004971127078137211, 0.004971127078137211, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 3.778782064908889, 34.00903858418, 0.001852364854527473, 0.001852364854527472)
  "ECs|M2"    = @(1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.004971127078137211, 0.004971127078137211, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 2.370243307037778, 21.33218976334, 0.001161896961300839, 0.001161896961300839)
  "ECs|sCs"   = @(1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.004971127078137211, 0.004971127078137211, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 1.074290393392778, 9.668613540535, 0.0005266188242918039, 0.0005266188242918039)
  "FAPs|ECs"  = @(3, 1, 14.289121, 42.867363, 0.760294375288143, 0.760294375288143, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 446.2344548267246, 4016.110093440522, 0.2187448248673166, 0.2187448248673166)
  "FAPs|FAPs" = @(3, 1, 14.289121, 42.867363, 0.760294375288143, 0.760294375288143, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 577.9346824637026, 5201.412142173323, 0.2833044816079039, 0.2833044816079039)
  "FAPs|M2"   = @(3, 1, 14.289121, 42.867363, 0.760294375288143, 0.760294375288143, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 362.5098747386014, 3262.588872647412, 0.1777029052881176, 0.1777029052881176)
  "FAPs|sCs"  = @(3, 1, 14.289121, 42.867363, 0.760294375288143, 0.760294375288143, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 164.3041770375903, 1478.737593338313, 0.08054216352480503, 0.08054216352480503)
  "sCs|ECs"   = @(3, 1, 4.411646, 13.234938, 0.2347344976337197, 0.2347344976337198, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 137.7711370558413, 1239.940233502572, 0.06753562599453092, 0.06753562599453093)
  "sCs|FAPs"  = @(3, 1, 4.411646, 13.234938, 0.2347344976337197, 0.2347344976337198, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 178.4324753182693, 1605.892277864424, 0.08746787735001915, 0.08746787735001917)
  "sCs|M2"    = @(3, 1, 4.411646, 13.234938, 0.2347344976337197, 0.2347344976337198, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 111.9218767096347, 1007.296890386712, 0.05486427830674137, 0.05486427830674138)
  "sCs|sCs"   = @(3, 1, 4.411646, 13.234938, 0.2347344976337197, 0.2347344976337198, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 50.72753358384868, 456.5478022546381, 0.02486671598242831, 0.02486671598242831)
}

$r = 2
foreach ($sender in $senders) {
  foreach ($target in $targets) {
    $ws.Cells.Item($r, 1).Value = $sender
    $ws.Cells.Item($r, 2).Value = "Angpt1"
    $ws.Cells.Item($r, 3).Value = "Itga5"
    $ws.Cells.Item($r, 4).Value = $target

    $row = $values["$sender|$target"]
    for ($i = 0; $i -lt $row.Length; $i++) {
      $ws.Cells.Item($r, 5 + $i).Value = $row[$i]
    }
    $r++
  }
}
